$d = $word.ActiveDocument

# Change 1: remove the word "nhóm " right before "chúng tôi quyết định thiết kế hệ thống với"
$d.Content.Find.Execute(
    "n, nhóm chúng tôi",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "n, chúng tôi",
    2
)
